Write-Host "START"
$p = $ppt.ActivePresentation
try {
  $p.ThisMethodTotallyDoesNotExist123()
  Write-Host "after call, no exception"
} catch {
  Write-Host "ERR1:" $_
}
Write-Host "END"
